$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Pre-format D2:D51 and E2:E51 as Text so numeric-looking strings
# (e.g. "1.003", "325.74") are stored as text, matching the source data,
# instead of being auto-converted to floating point numbers by Excel.
$ws.Range("D2:D51").NumberFormat = "@"
$ws.Range("E2:E51").NumberFormat = "@"

$ws.Range('D2').Value = '29.140.57'
$ws.Range('E2').Value = '  -0.46%  '
$ws.Range('D3').Value = '1.905.94'
$ws.Range('E3').Value = '  -0.47%  '
$ws.Range('D4').Value = '1.003'
$ws.Range('E4').Value = '  -0.04%  '
$ws.Range('D5').Value = '325.74'
$ws.Range('E5').Value = '  -0.64%  '
$ws.Range('D6').Value = '1.001'
$ws.Range('E6').Value = '  -0.15%  '
$ws.Range('D7').Value = '0.4619'
$ws.Range('E7').Value = '  -0.10%  '
$ws.Range('E8').Value = '  -1.45%  '
$ws.Range('D9').Value = '0.07872'
$ws.Range('E9').Value = '  -0.98%  '
$ws.Range('D10').Value = '0.9909'
$ws.Range('E10').Value = '  -1.26%  '
$ws.Range('D11').Value = '21.96'
$ws.Range('E11').Value = '  -2.02%  '
$ws.Range('D12').Value = '1.889.07'
$ws.Range('E12').Value = '  -1.50%  '
$ws.Range('D13').Value = '5.766'
$ws.Range('E13').Value = '  -0.09%  '
$ws.Range('D14').Value = '7.044'
$ws.Range('E14').Value = '  -0.96%  '
$ws.Range('D15').Value = '0.07049'
$ws.Range('E15').Value = '  +1.36%  '
$ws.Range('D16').Value = '88.11'
$ws.Range('E16').Value = '  -0.61%  '
$ws.Range('D17').Value = '1.004'
$ws.Range('E17').Value = '  +0.05%  '
$ws.Range('D18').Value = '0.000009921'
$ws.Range('E18').Value = '  -1.92%  '
$ws.Range('D19').Value = '17.08'
$ws.Range('E19').Value = '  -0.31%  '
$ws.Range('D20').Value = '1.001'
$ws.Range('E20').Value = '  +0.00%  '
$ws.Range('D21').Value = '29.159.01'
$ws.Range('E21').Value = '  -0.45%  '
$ws.Range('D22').Value = '5.327'
$ws.Range('E22').Value = '  -0.74%  '
$ws.Range('E23').Value = '  +0.01%  '
$ws.Range('D24').Value = '2.099'
$ws.Range('E24').Value = '  +1.67%  '
$ws.Range('D25').Value = '156.49'
$ws.Range('E25').Value = '  -0.22%  '
$ws.Range('D26').Value = '19.41'
$ws.Range('E26').Value = '  -0.61%  '
$ws.Range('D27').Value = '5.896'
$ws.Range('E27').Value = '  -3.39%  '
$ws.Range('D28').Value = '118.73'
$ws.Range('E28').Value = '  -0.33%  '
$ws.Range('D29').Value = '1.880'
$ws.Range('E29').Value = '  -6.09%  '
$ws.Range('D30').Value = '0.09355'
$ws.Range('E30').Value = '  -0.51%  '
$ws.Range('D31').Value = '0.8962'
$ws.Range('E31').Value = '  -3.54%  '
$ws.Range('D32').Value = '5.235'
$ws.Range('E32').Value = '  -2.33%  '
$ws.Range('D33').Value = '1.321'
$ws.Range('D34').Value = '3.148'
$ws.Range('E34').Value = '  -3.91%  '
$ws.Range('D35').Value = '0.05794'
$ws.Range('E35').Value = '  -0.90%  '
$ws.Range('D36').Value = '1.172'
$ws.Range('E36').Value = '  -3.01%  '
$ws.Range('D37').Value = '0.02089'
$ws.Range('E37').Value = '  -0.95%  '
$ws.Range('D38').Value = '1.000'
$ws.Range('E38').Value = '  -0.16%  '
$ws.Range('D39').Value = '0.5705'
$ws.Range('E39').Value = '  -1.04%  '
$ws.Range('D40').Value = '7.674'
$ws.Range('E40').Value = '  -3.83%  '
$ws.Range('D41').Value = '0.1813'
$ws.Range('E41').Value = '  +0.33%  '
$ws.Range('E42').Value = '  -2.92%  '
$ws.Range('D43').Value = '11.88'
$ws.Range('E43').Value = '  -1.42%  '
$ws.Range('D44').Value = '0.5358'
$ws.Range('E44').Value = '  -1.36%  '
$ws.Range('D45').Value = '2.176'
$ws.Range('E45').Value = '  -4.71%  '
$ws.Range('D46').Value = '0.07011'
$ws.Range('E46').Value = '  -0.96%  '
$ws.Range('E47').Value = '  -2.13%  '
$ws.Range('E48').Value = '  -0.62%  '
$ws.Range('D49').Value = '113.18'
$ws.Range('E49').Value = '  -0.21%  '
$ws.Range('D50').Value = '0.2998'
$ws.Range('E50').Value = '  +1.68%  '
$ws.Range('D51').Value = '71.31'
$ws.Range('E51').Value = '  -0.73%  '

# Remove the temporary Text number format again so the cells end up with
# the same (default/general) style they started with, matching the original
# workbook formatting -- only cell values should differ.
$ws.Range("D2:D51").ClearFormats()
$ws.Range("E2:E51").ClearFormats()
